$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# Rows 6 and 7 were blank placeholder rows (height-only). Insert two rows at
# position 6 so the new rows inherit row 5's formatting/style (s="2"), then
# remove the two rows that got pushed past the end (1001/1002) to keep the
# sheet at 1000 rows, matching the original layout.
$ws.Range("A6:A7").EntireRow.Insert()

# Row 6: OutputReportFolder (shared-string insertion order: B, A, C)
$ws.Range("B6").Value = "ADYF.OutputReportFolder"
$ws.Range("A6").Value = "OutputReportFolder"
$ws.Range("C6").Value = "Report Folder where files are downloaded"

# Row 7: ConsolidatedOutputFile (shared-string insertion order: B, A, C)
$ws.Range("B7").Value = "ADYF.OutputExcelFileName"
$ws.Range("A7").Value = "ConsolidatedOutputFile"
$ws.Range("C7").Value = "File to store the consolidated data from Downloaded reports"

# Undo the downward shift past row 1000 caused by the row insert above.
$ws.Range("A1001:A1002").EntireRow.Delete()

# Restore the row height metadata the insert operation dropped.
$ws.Rows(6).RowHeight = 14.25
$ws.Rows(7).RowHeight = 14.25

$ws.Activate()
$ws.Range("A7").Select()
